$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" date placeholder text from
#    06/23/2022 to 06/29/2022 on the slide master and every slide layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "06/23/2022") {
                $tr.Text = "06/29/2022"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 2) Fix the misspelled "Resource Mananger" -> "Resource Manager" text
#    on slide 12 (three shapes, each holding the text split across two
#    runs: "Resource " + "Mananger"). Merge them into a single run
#    that keeps the first run's formatting.
# ---------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
for ($i = 1; $i -le $slide12.Shapes.Count; $i++) {
    $shp = $slide12.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Resource Mananger") {
            # Delete the trailing "Mananger" run entirely...
            $tail = $tr.Characters(10, 8)
            $tail.Delete()
            # ...then append "Manager" inside the remaining "Resource "
            # run so it keeps that run's formatting as a single run.
            $head = $tr.Characters(1, 9)
            $head.InsertAfter("Manager") | Out-Null
        }
    }
}
